# Update values produced by an updated RandomForest imputation run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value  = 4.748900000000003
$ws.Range("D4").Value  = -8.2605
$ws.Range("D5").Value  = -8.696399999999997
$ws.Range("B6").Value  = 9.454000000000001
$ws.Range("B7").Value  = 6.510799999999997
$ws.Range("D8").Value  = -8.409899999999999
$ws.Range("B16").Value = 8.626300000000008
$ws.Range("D16").Value = -8.3622
$ws.Range("B20").Value = 5.638599999999997
$ws.Range("D22").Value = -8.073899999999998
